$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores figures as text (e.g. "25.387.40", "1.000")
# in the source workbook. Mark the cells whose refreshed values look like
# plain numbers as Text first so Excel keeps storing them as strings,
# matching the existing inline-string cells in this column.
$ws.Range('D4:D10').NumberFormat = '@'
$ws.Range('D12:D16').NumberFormat = '@'
$ws.Range('D19:D19').NumberFormat = '@'
$ws.Range('D22:D29').NumberFormat = '@'
$ws.Range('D31:D41').NumberFormat = '@'
$ws.Range('D43:D49').NumberFormat = '@'
$ws.Range('D51:D51').NumberFormat = '@'

# Apply the updated coin data (price + volume figures refreshed by the
# scheduled GitHub Actions run).
$ws.Range('D2').Value = '25.387.40'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '1.662.90'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').Value = '235.81'
$ws.Range('E5').Value = '  -2.11%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '0.4782'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('D8').Value = '0.2605'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('D9').Value = '0.06150'
$ws.Range('E9').Value = '  +1.43%  '
$ws.Range('D10').Value = '0.07070'
$ws.Range('E10').Value = '  -1.38%  '
$ws.Range('D11').Value = '1.662.01'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').Value = '14.73'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').Value = '0.5909'
$ws.Range('E13').Value = '  -6.73%  '
$ws.Range('D14').Value = '4.382'
$ws.Range('E14').Value = '  -6.13%  '
$ws.Range('D15').Value = '74.36'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = '25.385.08'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = '0.000006762'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = '1.873.39'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').Value = '4.441'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').Value = '8.641'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').Value = '5.330'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '133.41'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = '15.04'
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').Value = '1.403'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').Value = '104.18'
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('D29').Value = '1.686'
$ws.Range('E29').Value = '  -3.16%  '
$ws.Range('E30').Value = '  +3.44%  '
$ws.Range('D31').Value = '3.618'
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').Value = '0.07644'
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('D33').Value = '0.04372'
$ws.Range('E33').Value = '  -5.94%  '
$ws.Range('D34').Value = '0.9996'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = '2.606'
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').Value = '0.6117'
$ws.Range('E36').Value = '  +3.45%  '
$ws.Range('D37').Value = '0.9433'
$ws.Range('E37').Value = '  -2.72%  '
$ws.Range('D38').Value = '2.613'
$ws.Range('E38').Value = '  -2.54%  '
$ws.Range('D39').Value = '0.8535'
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').Value = '0.9999'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').Value = '0.01504'
$ws.Range('E41').Value = '  -4.45%  '
$ws.Range('E42').Value = '  -3.43%  '
$ws.Range('D43').Value = '98.09'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('D44').Value = '0.3764'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').Value = '4.662'
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.204'
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1109'
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('D48').Value = '0.05250'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').Value = '29.51'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.327'
$ws.Range('E51').Value = '  -2.11%  '
